$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and date range) ---
$ws.Range("A8").Value = "Volume 29   Number  47"
$ws.Range("C9").Value = "Report Covering the Week  11/21/2022  Through  11/27/2022"

# --- Type/style transition cells (text <-> numeric) ---
# Row 15: Rape - This Week (C), Last Week (G), %Chg (H) become blank-marker text cells
$ws.Cells.Item(15, 3).Value = "'0"
$ws.Cells.Item(15, 4).Copy()
$ws.Cells.Item(15, 3).PasteSpecial(-4122)

$ws.Cells.Item(15, 7).Value = "'0"
$ws.Cells.Item(15, 4).Copy()
$ws.Cells.Item(15, 7).PasteSpecial(-4122)

$ws.Cells.Item(15, 8).Value = "'***.*"
$ws.Cells.Item(15, 5).Copy()
$ws.Cells.Item(15, 8).PasteSpecial(-4122)

# Row 20: G.L.A. - This Week (C) becomes numeric
$ws.Cells.Item(20, 3).Value = 1
$ws.Cells.Item(20, 4).Copy()
$ws.Cells.Item(20, 3).PasteSpecial(-4122)

# Row 22: Transit - Last Week (D) and %Chg (E) become numeric
$ws.Cells.Item(22, 4).Value = 1
$ws.Cells.Item(22, 6).Copy()
$ws.Cells.Item(22, 4).PasteSpecial(-4122)

$ws.Cells.Item(22, 5).Value = -100
$ws.Cells.Item(23, 8).Copy()
$ws.Cells.Item(22, 5).PasteSpecial(-4122)

# Row 23: Housing - This Week (C) becomes blank-marker text cell
$ws.Cells.Item(23, 3).Value = "'0"
$ws.Cells.Item(27, 4).Copy()
$ws.Cells.Item(23, 3).PasteSpecial(-4122)

# Row 26: UCR Rape* - This Week (C) becomes blank-marker text cell
$ws.Cells.Item(26, 3).Value = "'0"
$ws.Cells.Item(27, 4).Copy()
$ws.Cells.Item(26, 3).PasteSpecial(-4122)

# --- Plain numeric value updates (rows 15-29) ---
$ws.Cells.Item(15, 13).Value = -16.666666666666
$ws.Cells.Item(15, 14).Value = -76.190476190476
$ws.Cells.Item(16, 4).Value = 4
$ws.Cells.Item(16, 5).Value = 0
$ws.Cells.Item(16, 7).Value = 11
$ws.Cells.Item(16, 8).Value = 0
$ws.Cells.Item(16, 9).Value = 159
$ws.Cells.Item(16, 10).Value = 156
$ws.Cells.Item(16, 11).Value = 1.923076923076
$ws.Cells.Item(16, 12).Value = 21.374045801526
$ws.Cells.Item(16, 13).Value = -26.388888888888
$ws.Cells.Item(16, 14).Value = -79.822335025380
$ws.Cells.Item(17, 3).Value = 6
$ws.Cells.Item(17, 4).Value = 3
$ws.Cells.Item(17, 5).Value = 100
$ws.Cells.Item(17, 7).Value = 19
$ws.Cells.Item(17, 8).Value = -21.052631578947
$ws.Cells.Item(17, 9).Value = 209
$ws.Cells.Item(17, 10).Value = 236
$ws.Cells.Item(17, 11).Value = -11.440677966101
$ws.Cells.Item(17, 12).Value = 5.025125628140
$ws.Cells.Item(17, 13).Value = 29.012345679012
$ws.Cells.Item(17, 14).Value = -68.991097922848
$ws.Cells.Item(18, 3).Value = 7
$ws.Cells.Item(18, 4).Value = 3
$ws.Cells.Item(18, 5).Value = 133.333333333333
$ws.Cells.Item(18, 6).Value = 10
$ws.Cells.Item(18, 7).Value = 12
$ws.Cells.Item(18, 8).Value = -16.666666666666
$ws.Cells.Item(18, 9).Value = 159
$ws.Cells.Item(18, 10).Value = 147
$ws.Cells.Item(18, 11).Value = 8.163265306122
$ws.Cells.Item(18, 12).Value = 17.777777777777
$ws.Cells.Item(18, 13).Value = 60.606060606060
$ws.Cells.Item(18, 14).Value = -80.174563591022
$ws.Cells.Item(19, 3).Value = 6
$ws.Cells.Item(19, 4).Value = 5
$ws.Cells.Item(19, 5).Value = 20
$ws.Cells.Item(19, 7).Value = 27
$ws.Cells.Item(19, 8).Value = -11.111111111111
$ws.Cells.Item(19, 9).Value = 361
$ws.Cells.Item(19, 10).Value = 349
$ws.Cells.Item(19, 11).Value = 3.438395415472
$ws.Cells.Item(19, 12).Value = 24.482758620689
$ws.Cells.Item(19, 13).Value = 25.347222222222
$ws.Cells.Item(19, 14).Value = -18.140589569161
$ws.Cells.Item(20, 4).Value = 1
$ws.Cells.Item(20, 5).Value = 0
$ws.Cells.Item(20, 6).Value = 4
$ws.Cells.Item(20, 7).Value = 6
$ws.Cells.Item(20, 8).Value = -33.333333333333
$ws.Cells.Item(20, 9).Value = 47
$ws.Cells.Item(20, 10).Value = 48
$ws.Cells.Item(20, 11).Value = -2.083333333333
$ws.Cells.Item(20, 12).Value = 9.302325581395
$ws.Cells.Item(20, 13).Value = 161.111111111111
$ws.Cells.Item(20, 14).Value = -64.925373134328
$ws.Cells.Item(21, 3).Value = 24
$ws.Cells.Item(21, 4).Value = 16
$ws.Cells.Item(21, 5).Value = 50
$ws.Cells.Item(21, 6).Value = 65
$ws.Cells.Item(21, 7).Value = 75
$ws.Cells.Item(21, 8).Value = -13.333333333333
$ws.Cells.Item(21, 9).Value = 948
$ws.Cells.Item(21, 10).Value = 945
$ws.Cells.Item(21, 11).Value = 0.317460317460
$ws.Cells.Item(21, 12).Value = 15.750915750915
$ws.Cells.Item(21, 13).Value = 18.5
$ws.Cells.Item(21, 14).Value = -67.456230690010
$ws.Cells.Item(22, 7).Value = 4
$ws.Cells.Item(22, 8).Value = -75
$ws.Cells.Item(22, 10).Value = 25
$ws.Cells.Item(22, 11).Value = 20
$ws.Cells.Item(23, 5).Value = -100
$ws.Cells.Item(23, 6).Value = 7
$ws.Cells.Item(23, 8).Value = 40
$ws.Cells.Item(23, 10).Value = 72
$ws.Cells.Item(23, 11).Value = 16.666666666666
$ws.Cells.Item(23, 13).Value = 90.909090909090
$ws.Cells.Item(24, 3).Value = 22
$ws.Cells.Item(24, 4).Value = 14
$ws.Cells.Item(24, 5).Value = 57.142857142857
$ws.Cells.Item(24, 6).Value = 96
$ws.Cells.Item(24, 7).Value = 67
$ws.Cells.Item(24, 8).Value = 43.283582089552
$ws.Cells.Item(24, 9).Value = 1059
$ws.Cells.Item(24, 10).Value = 1070
$ws.Cells.Item(24, 11).Value = -1.028037383177
$ws.Cells.Item(24, 12).Value = 31.064356435643
$ws.Cells.Item(24, 13).Value = 14.116379310344
$ws.Cells.Item(25, 3).Value = 7
$ws.Cells.Item(25, 4).Value = 2
$ws.Cells.Item(25, 5).Value = 250
$ws.Cells.Item(25, 6).Value = 29
$ws.Cells.Item(25, 7).Value = 28
$ws.Cells.Item(25, 8).Value = 3.571428571428
$ws.Cells.Item(25, 9).Value = 375
$ws.Cells.Item(25, 10).Value = 374
$ws.Cells.Item(25, 11).Value = 0.267379679144
$ws.Cells.Item(25, 12).Value = 12.951807228915
$ws.Cells.Item(25, 13).Value = -19.871794871794
$ws.Cells.Item(26, 7).Value = 1
$ws.Cells.Item(26, 8).Value = 200
$ws.Cells.Item(27, 6).Value = 4
$ws.Cells.Item(27, 7).Value = 6
$ws.Cells.Item(27, 8).Value = -33.333333333333
$ws.Cells.Item(28, 7).Value = 3
$ws.Cells.Item(28, 8).Value = -66.666666666666
$ws.Cells.Item(28, 10).Value = 28
$ws.Cells.Item(28, 11).Value = -14.285714285714
$ws.Cells.Item(29, 7).Value = 3
$ws.Cells.Item(29, 8).Value = -66.666666666666
$ws.Cells.Item(29, 10).Value = 25
$ws.Cells.Item(29, 11).Value = -24
